$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch remain stored as text (matches original inline-string cells),
# preventing Excel from auto-converting numeric-looking strings into numbers and
# stripping formatting such as trailing zeros.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.903.73"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "2.305.74"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "322.71"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "105.20"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "40.32"
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("D11").Value = "0.0911"
$ws.Range("D12").Value = "8.60"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "0.975"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "15.36"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "2.655.93"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "2.301.79"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "42.800.70"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  +35.10%  "
$ws.Range("D22").Value = "73.80"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "271.51"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "22.72"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "38.25"
$ws.Range("E30").Value = "  +11.19%  "
$ws.Range("D31").Value = "165.74"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  +6.24%  "
$ws.Range("D33").Value = "0.0887"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").Value = "0.133"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "2.55"
$ws.Range("E35").Value = "  -12.10%  "
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").Value = "1.56"
$ws.Range("E41").Value = "  +7.37%  "
$ws.Range("D42").Value = "99.18"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("D43").Value = "70.61"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.225"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "12.48"
$ws.Range("E46").Value = "  +5.56%  "
$ws.Range("D47").Value = "82.68"
$ws.Range("E47").Value = "  +9.65%  "
$ws.Range("D48").Value = "113.93"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "5.31"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.91"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "1.595.38"
$ws.Range("E51").Value = "  +4.66%  "
